$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 48 (row 57): RCB vs PBKS
$ws.Range("E57").Value = 100
$ws.Range("H57").Value = 60
$ws.Range("K57").Value = 40
$ws.Range("N57").Value = 80
$ws.Range("Q57").Value = 0
$ws.Range("T57").Value = 20

# Contest 49 (row 58): KKR vs SRH
$ws.Range("E58").Value = 80
$ws.Range("H58").Value = 40
$ws.Range("K58").Value = 20
$ws.Range("N58").Value = 60
$ws.Range("Q58").Value = 0
$ws.Range("T58").Value = 100

# Scroll the frozen pane's top-left visible cell down two rows (A42 -> A44)
$ws.Application.ActiveWindow.ScrollRow = 44
